# CastelfrancoE.xlsx update: add 2 new days of data (rows), shifting the
# existing "ultimi 7 giorni" rolling-window data down by one row and
# recomputing the rolling sums/rates for the affected window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 93 (pushes old rows 93-113 to 94-114) ---
$ws.Rows(93).Insert()

# Fix formatting of the newly inserted row 93 (copy format from row 94,
# which holds the same style the row had before the insert shifted it).
$ws.Range("A94:D94").Copy()
$ws.Range("A93:D93").PasteSpecial(-4122)

# --- 2. Append a brand-new row 115 after the current last row (114) ---
$ws.Range("A114:D114").Copy()
$ws.Range("A115:D115").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 3. Write the date (col A) / new-positives (col B) values ---
$ws.Cells.Item(93,1).Value = 44235
$ws.Cells.Item(93,2).Value = 10

$ws.Cells.Item(115,1).Value = 44257
$ws.Cells.Item(115,2).Value = 24

# --- 4. Recompute the rolling 7-day sum (col C) and rate per 100k (col D)
#         for every row whose centered 7-row window changed. Values taken
#         from the source data so the stored doubles match bit-for-bit.
$ws.Cells.Item(90,3).Value = 69
$ws.Cells.Item(90,4).Value = 209.1986781069035

$ws.Cells.Item(91,3).Value = 77
$ws.Cells.Item(91,4).Value = 233.4535973077039

$ws.Cells.Item(92,3).Value = 75
$ws.Cells.Item(92,4).Value = 227.3898675075039

$ws.Cells.Item(93,3).Value = 67
$ws.Cells.Item(93,4).Value = 203.1349483067035

$ws.Cells.Item(94,3).Value = 73
$ws.Cells.Item(94,4).Value = 221.3261377073038

$ws.Cells.Item(95,3).Value = 79
$ws.Cells.Item(95,4).Value = 239.5173271079041

$ws.Cells.Item(96,3).Value = 79
$ws.Cells.Item(96,4).Value = 239.5173271079041

$ws.Cells.Item(97,3).Value = 93
$ws.Cells.Item(97,4).Value = 281.9634357093048

$ws.Cells.Item(98,3).Value = 103
$ws.Cells.Item(98,4).Value = 312.2820847103053

$ws.Cells.Item(99,3).Value = 102
$ws.Cells.Item(99,4).Value = 309.2502198102052

$ws.Cells.Item(100,3).Value = 108
$ws.Cells.Item(100,4).Value = 327.4414092108056

$ws.Cells.Item(101,3).Value = 120
$ws.Cells.Item(101,4).Value = 363.8237880120062

$ws.Cells.Item(102,3).Value = 124
$ws.Cells.Item(102,4).Value = 375.9512476124064

$ws.Cells.Item(103,3).Value = 139
$ws.Cells.Item(103,4).Value = 421.4292211139072

$ws.Cells.Item(104,3).Value = 130
$ws.Cells.Item(104,4).Value = 394.1424370130067

$ws.Cells.Item(105,3).Value = 126
$ws.Cells.Item(105,4).Value = 382.0149774126065

$ws.Cells.Item(106,3).Value = 128
$ws.Cells.Item(106,4).Value = 388.0787072128066

$ws.Cells.Item(107,3).Value = 129
$ws.Cells.Item(107,4).Value = 391.1105721129067

$ws.Cells.Item(108,3).Value = 132
$ws.Cells.Item(108,4).Value = 400.2061668132068

$ws.Cells.Item(109,3).Value = 142
$ws.Cells.Item(109,4).Value = 430.5248158142073

$ws.Cells.Item(110,3).Value = 142
$ws.Cells.Item(110,4).Value = 430.5248158142073

$ws.Cells.Item(111,3).Value = 155
$ws.Cells.Item(111,4).Value = 469.939059515508

$ws.Cells.Item(112,3).Value = 159
$ws.Cells.Item(112,4).Value = 482.0665191159081

# Rows 113, 114, 115 don't have a full 7-row window yet (no future data),
# so their C/D stay blank - same as the source rows they were copied from.

Write-Output "CastelfrancoE.xlsx update applied"
